$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (MSE) and C (MAE) hold numeric-looking values stored as text.
# Force text format before writing so they keep their exact string form.
$textRange = $ws.Range("B2:C5")
$textRange.NumberFormat = "@"

# Row 2 (IIT VAL)
$ws.Range("B2").Value = "0.06565217"
$ws.Range("C2").Value = "0.12641451"

# Row 3 (REG VAL)
$ws.Range("B3").Value = "0.07652219"
$ws.Range("C3").Value = "0.13153398"
$ws.Range("D3").Value = 74
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4 (IIT TEST)
$ws.Range("B4").Value = "0.1771727"
$ws.Range("C4").Value = "0.3283748"

# Row 5 (REG TEST)
$ws.Range("B5").Value = "0.1434832"
$ws.Range("C5").Value = "0.2826047"
